$d = $word.ActiveDocument

# Locate the paragraph that currently reads "Please avoid me."
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Please avoid me.") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    $target = $d.Paragraphs.Item(2)
}

$full = $target.Range
$start = $full.Start
$end = $full.End - 1   # exclude the paragraph mark

# Replace the run's text with the trailing fragment "avoid me." first. Doing this
# in one shot (rather than clearing to "") keeps the paragraph non-empty at every
# step, which keeps later Range offsets valid.
$tail = $d.Range($start, $end)
$tail.Text = "avoid me."

# Insert "don't " (curly apostrophe, matching Word's AutoCorrect) right before it,
# as its own run.
$apostrophe = [char]0x2019
$dontRange = $d.Range($tail.Start, $tail.Start)
$dontRange.InsertBefore("don" + $apostrophe + "t ")

# The "_GoBack" bookmark (currently wrapping the picture paragraph, courtesy of
# Word's "last edit" tracking) now belongs between "don't " and "avoid me." --
# re-adding it under the same reserved name moves it here.
$bookmarkPoint = $d.Range($dontRange.End, $dontRange.End)
$d.Bookmarks.Add("_GoBack", $bookmarkPoint)

# Finally insert "Please " as its own leading run.
$pleaseRange = $d.Range($dontRange.Start, $dontRange.Start)
$pleaseRange.InsertBefore("Please ")
